$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.078.30'
$ws.Range('E2').Value = '  -0.90%  '
$ws.Range('D3').Value = '1.648.75'
$ws.Range('E3').Value = '  -0.98%  '
$ws.Range('E4').Value = '  -0.47%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.10'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5183'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.99%  '
$ws.Range('E7').Value = '  -0.44%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2610'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.89%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06274'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.99%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.44'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07798'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.47%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.454'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.47%  '
$ws.Range('D13').Value = '1.671.71'
$ws.Range('E13').Value = '  +0.18%  '
$ws.Range('D14').Value = '1.875.95'
$ws.Range('E14').Value = '  -0.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5531'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.09%  '
$ws.Range('D16').Value = '0.0₅7972'
$ws.Range('E16').Value = '  -3.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.65'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.73%  '
$ws.Range('D18').Value = '26.083.36'
$ws.Range('E18').Value = '  -0.95%  '
$ws.Range('E19').Value = '  -0.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.622'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '193.82'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.06'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.87%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.934'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.77%  '
$ws.Range('E24').Value = '  -0.41%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.93'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.43%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1203'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.31%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.165'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.42%  '
$ws.Range('E28').Value = '  -1.54%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.473'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.80%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05594'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.24%  '
$ws.Range('E31').Value = '  -1.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.469'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.32%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.377'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.91%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.592'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.54%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.800'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9463'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.403'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5643'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.90%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.958'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.46%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01575'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.89%  '
$ws.Range('D41').Value = '1.058.59'
$ws.Range('E41').Value = '  +0.57%  '
$ws.Range('E42').Value = '  -0.55%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8384'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.48%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '102.50'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.50%  '
$ws.Range('D45').Value = '1.788.06'
$ws.Range('E45').Value = '  -0.86%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.01'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E47').Value = '  +4.29%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05356'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.70%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.005'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.97%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4337'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.926'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.26%  '
